# Refresh the cryptos price/volume snapshot per the Fri Apr 21 03:26:33 UTC 2023
# GitHub Actions data pull. Sheet1 stores every data cell as text, so values
# that look numeric are written with a leading apostrophe (quote-prefix) to
# stop Excel from silently re-typing them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.565.52'
$ws.Range("E2").Value = '  -1.73%  '

$ws.Range("D3").Value = '1.968.82'
$ws.Range("E3").Value = '  +0.37%  '

$ws.Range("E4").Value = '  +0.27%  '

$ws.Range("D5").Value = '''323.47'
$ws.Range("E5").Value = '  -1.12%  '

$ws.Range("D6").Value = '''1.011'
$ws.Range("E6").Value = '  +0.21%  '

$ws.Range("D7").Value = '''0.4799'
$ws.Range("E7").Value = '  -3.52%  '

$ws.Range("D8").Value = '''0.4070'
$ws.Range("E8").Value = '  -3.13%  '

$ws.Range("D9").Value = '''54.11'
$ws.Range("E9").Value = '  +0.44%  '

$ws.Range("D10").Value = '''0.08525'
$ws.Range("E10").Value = '  -4.98%  '

$ws.Range("D11").Value = '''1.065'
$ws.Range("E11").Value = '  -2.88%  '

$ws.Range("D12").Value = '''22.54'
$ws.Range("E12").Value = '  -1.60%  '

$ws.Range("D13").Value = '2.055.49'
$ws.Range("E13").Value = '  +5.94%  '

$ws.Range("D14").Value = '''7.651'
$ws.Range("E14").Value = '  -2.57%  '

$ws.Range("D15").Value = '''6.211'
$ws.Range("E15").Value = '  -3.17%  '

$ws.Range("E16").Value = '  +0.08%  '

$ws.Range("D17").Value = '''91.48'
$ws.Range("E17").Value = '  +0.54%  '

$ws.Range("D18").Value = '''0.00001077'
$ws.Range("E18").Value = '  -1.58%  '

$ws.Range("D19").Value = '''0.06651'
$ws.Range("E19").Value = '  +0.33%  '

$ws.Range("D20").Value = '''18.68'
$ws.Range("E20").Value = '  -2.50%  '

$ws.Range("E21").Value = '  +0.24%  '

$ws.Range("D22").Value = '''5.889'
$ws.Range("E22").Value = '  -1.11%  '

$ws.Range("D23").Value = '28.596.09'
$ws.Range("E23").Value = '  -1.69%  '

$ws.Range("D24").Value = '''11.54'
$ws.Range("E24").Value = '  -3.32%  '

$ws.Range("D25").Value = '''2.300'
$ws.Range("E25").Value = '  +0.52%  '

$ws.Range("D26").Value = '2.225.89'
$ws.Range("E26").Value = '  +1.95%  '

$ws.Range("D27").Value = '''156.38'
$ws.Range("E27").Value = '  +0.26%  '

$ws.Range("D28").Value = '''20.41'
$ws.Range("E28").Value = '  -0.77%  '

$ws.Range("D29").Value = '''5.916'
$ws.Range("E29").Value = '  -4.71%  '

$ws.Range("D30").Value = '''2.187'
$ws.Range("E30").Value = '  -3.00%  '

$ws.Range("D31").Value = '''125.05'
$ws.Range("E31").Value = '  -1.37%  '

$ws.Range("D32").Value = '''0.9919'
$ws.Range("E32").Value = '  -4.62%  '

$ws.Range("D33").Value = '''0.09684'
$ws.Range("E33").Value = '  -1.39%  '

$ws.Range("D34").Value = '''1.467'
$ws.Range("E34").Value = '  -3.96%  '

$ws.Range("D35").Value = '''5.654'
$ws.Range("E35").Value = '  -2.42%  '

$ws.Range("D36").Value = '''3.699'
$ws.Range("E36").Value = '  +0.14%  '

$ws.Range("D37").Value = '''9.178'
$ws.Range("E37").Value = '  +2.52%  '

$ws.Range("D38").Value = '''0.02342'
$ws.Range("E38").Value = '  -3.06%  '

$ws.Range("D39").Value = '''0.06257'
$ws.Range("E39").Value = '  -0.76%  '

$ws.Range("D40").Value = '''1.261'
$ws.Range("E40").Value = '  -1.93%  '

$ws.Range("D41").Value = '''0.6248'
$ws.Range("E41").Value = '  -2.67%  '

$ws.Range("D42").Value = '''11.26'
$ws.Range("E42").Value = '  -1.26%  '

$ws.Range("E43").Value = '  +0.21%  '

$ws.Range("D44").Value = '''0.1926'
$ws.Range("E44").Value = '  -2.99%  '

$ws.Range("D45").Value = '''1.355'
$ws.Range("E45").Value = '  +5.56%  '

$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '''0.5969'
$ws.Range("E46").Value = '  -3.31%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '''13.07'
$ws.Range("E47").Value = '  -2.81%  '

$ws.Range("D48").Value = '''2.077'
$ws.Range("E48").Value = '  -4.24%  '

$ws.Range("D49").Value = '''3.419'

$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").Value = '''0.00000000317'
$ws.Range("E50").Value = '  -2.87%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '''0.06824'
$ws.Range("E51").Value = '  -0.65%  '
